$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "580.17") must be
# forced to stay text-typed (matching the original inlineStr cells), since
# Excel auto-converts plain-looking numeric strings to numbers otherwise.
# Pre-format as text, assign, then restore the Normal style so no stray
# number-format is left applied to the cell.
$numericLookingCells = @("D4","D5","D6","D8","D10","D12","D13","D15","D19","D20","D21","D22","D24","D25","D27","D28","D31","D32","D34","D36","D37","D38","D39","D42","D43","D44","D45","D47","D49","D51")
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '66.937.22'
$ws.Range("E2").Value = '  -1.25%  '
$ws.Range("D3").Value = '2.459.30'
$ws.Range("E3").Value = '  -1.58%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '580.17'
$ws.Range("E5").Value = '  -2.06%  '
$ws.Range("D6").Value = '166.21'
$ws.Range("E6").Value = '  -4.13%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '0.513'
$ws.Range("E8").Value = '  -2.57%  '
$ws.Range("D9").Value = '2.458.73'
$ws.Range("E9").Value = '  -1.52%  '
$ws.Range("D10").Value = '0.134'
$ws.Range("E10").Value = '  -4.12%  '
$ws.Range("E11").Value = '  -1.06%  '
$ws.Range("D12").Value = '4.89'
$ws.Range("E12").Value = '  -3.80%  '
$ws.Range("D13").Value = '0.332'
$ws.Range("E13").Value = '  -3.03%  '
$ws.Range("D14").Value = '2.903.19'
$ws.Range("E14").Value = '  -1.74%  '
$ws.Range("D15").Value = '25.32'
$ws.Range("E15").Value = '  -3.70%  '
$ws.Range("D16").Value = '66.592.35'
$ws.Range("E16").Value = '  -1.68%  '
$ws.Range("E17").Value = '  -4.91%  '
$ws.Range("D18").Value = '2.445.07'
$ws.Range("E18").Value = '  -1.92%  '
$ws.Range("D19").Value = '11.34'
$ws.Range("E19").Value = '  -3.80%  '
$ws.Range("D20").Value = '7.65'
$ws.Range("E20").Value = '  -3.98%  '
$ws.Range("D21").Value = '353.57'
$ws.Range("E21").Value = '  -2.68%  '
$ws.Range("D22").Value = '4.01'
$ws.Range("E22").Value = '  -2.94%  '
$ws.Range("D24").Value = '69.18'
$ws.Range("E24").Value = '  -2.83%  '
$ws.Range("D25").Value = '4.21'
$ws.Range("E25").Value = '  -7.71%  '
$ws.Range("E26").Value = '  -8.07%  '
$ws.Range("D27").Value = '8.87'
$ws.Range("E27").Value = '  -9.77%  '
$ws.Range("D28").Value = '0.998'
$ws.Range("E28").Value = '  -0.23%  '
$ws.Range("D29").Value = '2.579.97'
$ws.Range("E29").Value = '  -1.56%  '
$ws.Range("D30").Value = '0.0₃0894'
$ws.Range("E30").Value = '  -7.55%  '
$ws.Range("D31").Value = '506.59'
$ws.Range("E31").Value = '  -4.78%  '
$ws.Range("D32").Value = '7.75'
$ws.Range("E32").Value = '  -5.97%  '
$ws.Range("E33").Value = '  -5.60%  '
$ws.Range("D34").Value = '1.22'
$ws.Range("E34").Value = '  -6.68%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").Value = '157.99'
$ws.Range("E36").Value = '  -0.42%  '
$ws.Range("D37").Value = '0.116'
$ws.Range("E37").Value = '  -9.02%  '
$ws.Range("D38").Value = '18.56'
$ws.Range("E38").Value = '  -0.49%  '
$ws.Range("D39").Value = '18.43'
$ws.Range("E39").Value = '  -0.94%  '
$ws.Range("E40").Value = '  -6.47%  '
$ws.Range("E41").Value = '  +0.17%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '1.66'
$ws.Range("E42").Value = '  -6.88%  '
$ws.Range("B43").Value = 'PolygonEcosystemToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D43").Value = '0.326'
$ws.Range("E43").Value = '  -6.38%  '
$ws.Range("D44").Value = '4.76'
$ws.Range("E44").Value = '  -7.01%  '
$ws.Range("D45").Value = '38.74'
$ws.Range("E45").Value = '  -2.72%  '
$ws.Range("E46").Value = '  -7.54%  '
$ws.Range("D47").Value = '140.92'
$ws.Range("E47").Value = '  -3.05%  '
$ws.Range("E48").Value = '  -6.08%  '
$ws.Range("D49").Value = '0.513'
$ws.Range("E49").Value = '  -6.53%  '
$ws.Range("D50").Value = '0.0₆0252'
$ws.Range("E50").Value = '  -7.41%  '
$ws.Range("B51").Value = 'Optimism'
$ws.Range("C51").Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range("D51").Value = '1.58'
$ws.Range("E51").Value = '  -7.14%  '

foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}
